# Prog-entera: Completa datos del problema 11.
# Fills in the Costo/Capacidad (Costos-Capacidades), Demanda (Demandas)
# columns, renames the "Costo (USD)" header to "Costo (USD/m2)", and adds a
# new "Otros-parametros" sheet with three extra model parameters.

$wb = $excel.ActiveWorkbook

$xlRight  = -4152
$xlNone   = -4142

# ---------------------------------------------------------------------
# 1) Costos-Capacidades: header rename + fill Costo (USD/m2) / Capacidad
# ---------------------------------------------------------------------
$wsCostos = $wb.Worksheets.Item("Costos-Capacidades")

$wsCostos.Range("C1").Value = "Costo (USD/m2)"

$costos = @(
    @(2,  533.82, 42),
    @(3,  180.23, 59),
    @(4,  455.65, 42),
    @(5,  464.15, 60),
    @(6,  557.34, 39),
    @(7,  515.08, 49),
    @(8,  286.82, 59),
    @(9,  244.48, 62),
    @(10, 497.22, 61),
    @(11, 572.32, 51),
    @(12, 162.46, 36),
    @(13, 249.07, 52)
)

foreach ($row in $costos) {
    $r = $row[0]

    $c = $wsCostos.Cells.Item($r, 3)
    $c.Value = $row[1]
    $c.NumberFormat = "General"
    $c.HorizontalAlignment = $xlRight
    $c.Borders.LineStyle = $xlNone

    $d = $wsCostos.Cells.Item($r, 4)
    $d.Value = $row[2]
    $d.NumberFormat = "General"
    $d.HorizontalAlignment = $xlRight
}

# ---------------------------------------------------------------------
# 2) Demandas: fill Demanda column
# ---------------------------------------------------------------------
$wsDemandas = $wb.Worksheets.Item("Demandas")

$demandas = @(
    @(2, 27), @(3, 30), @(4, 21), @(5, 27), @(6, 28), @(7, 18),
    @(8, 28), @(9, 27), @(10, 25), @(11, 29), @(12, 34), @(13, 17),
    @(14, 16), @(15, 29), @(16, 10), @(17, 24), @(18, 22), @(19, 23),
    @(20, 20), @(21, 6), @(22, 6), @(23, 17), @(24, 14), @(25, 13)
)

foreach ($row in $demandas) {
    $r = $row[0]
    $c = $wsDemandas.Cells.Item($r, 3)
    $c.Value = $row[1]
    $c.NumberFormat = "General"
    $c.HorizontalAlignment = $xlRight
}

# ---------------------------------------------------------------------
# 3) New sheet "Otros-parametros" at the end of the workbook
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsOtros = $wb.Worksheets.Add($null, $lastSheet)
$wsOtros.Name = "Otros-parametros"

$wsOtros.Cells.Item(1, 1).Value = "Parametro"
$wsOtros.Cells.Item(1, 2).Value = "Valor"
$wsOtros.Cells.Item(1, 3).Value = "Unidades"

$wsOtros.Cells.Item(2, 1).Value = "Area Fija"
$wsOtros.Cells.Item(2, 2).Value = 200
$wsOtros.Cells.Item(2, 3).Value = "m2"

$wsOtros.Cells.Item(3, 1).Value = "Area por Estacionamiento"
$wsOtros.Cells.Item(3, 2).Value = 15
$wsOtros.Cells.Item(3, 3).Value = "m2"

$wsOtros.Cells.Item(4, 1).Value = "Presupuesto"
$wsOtros.Cells.Item(4, 2).Value = 2400000
$wsOtros.Cells.Item(4, 3).Value = "USD"

# Re-use the same look-and-feel as the other data sheets (bold/bordered
# header, bordered body with the label column left-aligned and the
# value/unit columns centered) by cloning formats from "Origenes".
$wsOrigenes = $wb.Worksheets.Item("Origenes")

$wsOrigenes.Range("A1:C1").Copy()
$wsOtros.Range("A1:C1").PasteSpecial(-4122)

$wsOrigenes.Range("B2").Copy()
$wsOtros.Range("A2:A4").PasteSpecial(-4122)

$wsOrigenes.Range("A2").Copy()
$wsOtros.Range("B2:C4").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 4) Restore "Origenes" as the active/selected sheet (tab 0)
# ---------------------------------------------------------------------
$wb.Worksheets.Item("Origenes").Activate()
